# Add the new "2021年" row (row 12) to the bottom of the data table,
# mirroring the existing rows (2011年..2020年) in columns A:G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Year label in column A, formatted like the other year cells (A2:A11).
$ws.Range("A12").Value = "2021年"
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats

# Numeric data points for 2021.
$ws.Range("B12").Value = 47
$ws.Range("D12").Value = 785
$ws.Range("G12").Value = 738
